# removed thresholding and seems made the mask better
# Adds a new "actual roof" pixel-count column (B) for the existing houses,
# appends four new house rows (7-10) with full data, and highlights a few
# rows (the ones whose mask "got better") with a yellow fill. Row 10's
# recomputed RGB-mask value now wraps text in its cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new "actual roof" values for the already-present houses (col B) ---
$ws.Range("B3").Value = 2920
$ws.Range("B4").Value = 2894
$ws.Range("B5").Value = 3036
$ws.Range("B6").Value = 4358

# --- four brand-new house rows ---
$newRows = @(
    @{ Row = 7;  A = "ae7a28a6-02bc-4fa6-858b-ab1a6dc5a99f"; B = 2240; C = 13210.1518799;     D = 2817.6672399600002;  E = 2334.6312535799998 },
    @{ Row = 8;  A = "b718aabc-6b8a-42c9-920e-3378addd5810"; B = 2947; C = 2968.9877826299999; D = 2456.5923633000002;  E = 2298.8461954200002 },
    @{ Row = 9;  A = "b817e6c4-4176-4211-a5c8-77b54fe2e04e"; B = 2398; C = 6963.0494914399997; D = 5473.3906292499996;  E = 5382.5096341799999 },
    @{ Row = 10; A = "ba9b993d-c107-451f-9ba2-c6c8a79a18a2"; B = 2390; C = 1475.1629000299999; D = 1652.5141847699999;  E = 1475.1629000299999 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# --- wrap text for the recomputed RGB-mask value in row 10 ---
$ws.Range("C10").WrapText = $true

# --- highlight the rows whose mask improved with a yellow fill ---
$yellow = 65535  # RGB(255,255,0) packed as an OLE BGR color
$ws.Range("A5").Interior.Color = $yellow
$ws.Range("A9").Interior.Color = $yellow
$ws.Range("A10").Interior.Color = $yellow

# --- selection moved further down the (now longer) sheet ---
$ws.Range("E22").Select()
